$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column F: format/codec for the existing rows 1-6 ----
$ws.Range("F1").Value = "AVC100"
$ws.Range("F2").Value = "AVC101"
$ws.Range("F3").Value = "AVC102"
$ws.Range("F4").Value = "AVC103"
$ws.Range("F5").Value = "AVC104"
$ws.Range("F6").Value = "AVC105"

# ---- New rows 7-9: BT advertisement inserts ----
$ws.Range("A7").Value = "BTVP1007"
$ws.Range("B7").Value = "Ology"
$ws.Range("E7").Value = "00:00:50:22"
$ws.Range("E7").NumberFormat = "h:mm:ss"
$ws.Range("F7").Value = "XDCAM50"

$ws.Range("A8").Value = "BTVP1008"
$ws.Range("B8").Value = "BT AD"
$ws.Range("E8").Value = "00:02:23:00"
$ws.Range("F8").Value = "XDCAM50"

$ws.Range("A9").Value = "BTVP1009"
$ws.Range("B9").Value = "BT AD"
$ws.Range("E9").Value = "00:01:01:03"
$ws.Range("F9").Value = "XDCAM50"

# ---- New rows 10-18: additional asset IDs (column A only) ----
$ws.Range("A10").Value = "BTVP1010"
$ws.Range("A11").Value = "BTVP1011"
$ws.Range("A12").Value = "BTVP1012"
$ws.Range("A13").Value = "BTVP1013"
$ws.Range("A14").Value = "BTVP1014"
$ws.Range("A15").Value = "BTVP1015"
$ws.Range("A16").Value = "BTVP1016"
$ws.Range("A17").Value = "BTVP1017"
$ws.Range("A18").Value = "BTVP1018"

# ---- Restore the selection highlighted in the sheet view ----
$null = $ws.Range("E7:E9").Select()
